$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top, shifting existing data down
$ws.Rows("1:2").Insert()

# New header row (row 1)
$ws.Range("A1").Value = "factura"
$ws.Range("B1").Value = "cod_serie"
$ws.Range("C1").Value = "num_serie"
$ws.Range("D1").Value = "destino"

# New data row (row 2)
$ws.Range("A2").Value = 20601452651
$ws.Range("B2").Value = "E001"
$ws.Range("C2").Value = 12654
$ws.Range("D2").Value = "ALMACEN COCINA LA MOLINA"

# Update selection to match target
$ws.Range("K9").Select()
